$d = $word.ActiveDocument

function Get-ParaByExactText($text) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $t2 = $t.Substring(0, $t.Length - 1)
        } else {
            $t2 = $t
        }
        if ($t2 -eq $text) {
            return $p
        }
    }
    return $null
}

function Set-ParaBodyXml($para, $innerXml) {
    $rng = $para.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Change 1: split the " (at least in headers)" run into three runs with
# expanded wording, inside the "Remove dependency on Windows..." paragraph.
# ---------------------------------------------------------------------------
$pRemove = Get-ParaByExactText("Remove dependency on Windows header files if possible (at least in headers).")
if ($pRemove -eq $null) { throw "Could not find 'Remove dependency' paragraph" }

$removeInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Remove dependency on Windows </w:t></w:r>' + `
  '<w:r><w:t>header files if possible</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> (at least in </w:t></w:r>' + `
  '<w:r><w:t>library headers, implementation is still permitted to use windows headers</w:t></w:r>' + `
  '<w:r><w:t>)</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>'
Set-ParaBodyXml $pRemove $removeInner

# ---------------------------------------------------------------------------
# Change 2: insert two new paragraphs around "Simplify Call..." and move the
# _GoBack bookmark from "Perform correct overflow checking." to a new
# paragraph about detail/windows.h.
# ---------------------------------------------------------------------------
$pSimplify = Get-ParaByExactText([string]::Format("Simplify Call if appropriate (don{0}t be too clever).", [char]0x2019))
if ($pSimplify -eq $null) { throw "Could not find 'Simplify Call' paragraph" }
$pOverflow = Get-ParaByExactText("Perform correct overflow checking.")
if ($pOverflow -eq $null) { throw "Could not find 'Perform correct overflow checking' paragraph" }

$spanRng = $d.Range($pSimplify.Range.Start, $pOverflow.Range.End)

$quote1 = [char]0x2018
$quote2 = [char]0x2019
$autoText = "Use " + $quote1 + "auto" + $quote2 + " where appropriate."

$spanInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Use a new detail/windows.h header to hold typedefs (namespaced) and __dllimports.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Simplify Call if appropriate (don' + $quote2 + 't be too clever).</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>' + $autoText + '</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Perform correct overflow checking.</w:t></w:r>' + `
  '</w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $spanInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $spanRng.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change 3: move <w:lastRenderedPageBreak/> from "Cache base pointers..." to
# "Extra sanity checking in all components."
# ---------------------------------------------------------------------------
$pExtra = Get-ParaByExactText("Extra sanity checking in all components.")
if ($pExtra -eq $null) { throw "Could not find 'Extra sanity checking' paragraph" }
$extraInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:lastRenderedPageBreak/><w:t>Extra sanity checking in all components.</w:t></w:r>' + `
  '</w:p>'
Set-ParaBodyXml $pExtra $extraInner

$cacheText = "Cache base pointers etc rather than retrieving it manually in every getter/setter. Slightly less " + $quote1 + "robust" + $quote2 + ", but due to the typically " + $quote1 + "read-only" + $quote2 + " nature of the data this is the expected behaviour in all known cases anyway."
$pCache = Get-ParaByExactText($cacheText)
if ($pCache -eq $null) { throw "Could not find 'Cache base pointers' paragraph" }
$cacheInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>' + $cacheText + '</w:t></w:r>' + `
  '</w:p>'
Set-ParaBodyXml $pCache $cacheInner

# ---------------------------------------------------------------------------
# Change 4: move <w:lastRenderedPageBreak/> from "Unknown value scan." to
# "Support injected scanning."
# ---------------------------------------------------------------------------
$pSupport = Get-ParaByExactText("Support injected scanning.")
if ($pSupport -eq $null) { throw "Could not find 'Support injected scanning' paragraph" }
$supportInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:lastRenderedPageBreak/><w:t>Support injected scanning.</w:t></w:r>' + `
  '</w:p>'
Set-ParaBodyXml $pSupport $supportInner

$pUnknown = Get-ParaByExactText("Unknown value scan.")
if ($pUnknown -eq $null) { throw "Could not find 'Unknown value scan' paragraph" }
$unknownInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Unknown value scan.</w:t></w:r>' + `
  '</w:p>'
Set-ParaBodyXml $pUnknown $unknownInner

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
